# Update Product-Backlog 3 and Sprint-Backlog 3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Product Backlog Item #1 in A2:A3 merge) ---
$ws.Range("A2").Value = "As a / an website visitor I want to Thấy công dụng việc làm video. so that Biết lợi ích khi có một video giới thiệu về cá nhân."
$ws.Range("B2").Value = "Tìm hiểu về công dụng khi làm video đó (mục đích)"
$ws.Range("C2").Value = "Tạ Việt Tiến"
$ws.Range("E2").Value = 3

# --- Row 3 (continues A2:A3 merge) ---
$ws.Range("B3").Value = "Viết nội dung phần tìm hiểu của Tiến lên web"
$ws.Range("C3").Value = "Trần Quang Thắng"
$ws.Range("E3").Value = 3

# --- Row 4 (Product Backlog Item #2 in A4:A6 merge) ---
$ws.Range("A4").Value = "As a / an website visitor I want to Thấy được những điều cần lưu ý (nội dung, background, nhạc, công cụ, …) khi làm một video  so that Dễ gây ấn tượng cho người xem hơn."
$ws.Range("B4").Value = "Tìm hiểu những nội dung cần có khi làm video"
$ws.Range("C4").Value = "Nguyễn Đào Xuân Trường"
$ws.Range("E4").Value = 3

# --- Row 5 (continues A4:A6 merge) ---
$ws.Range("B5").Value = "Tìm hiểu những lưu ý về background, nhạc,"
$ws.Range("C5").Value = "Nguyễn Quốc Huy"
$ws.Range("E5").Value = 3

# --- Row 6 (continues A4:A6 merge) ---
$ws.Range("B6").Value = "Viết nội dung phần tìm hiểu các bạn lên web"
$ws.Range("E6").Value = 3

# --- Update active selection to C11 ---
$ws.Range("C11").Select()
